$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Update the artificial-image generation values (columns ux/uy/sx/sy/theta/A
#     for stades 2/3/4/5, rows 28-33 of Sheet1) ---

# row 28 - ux
$ws.Range("B28").Value = 233.33
$ws.Range("E28").Value = 166.67
$ws.Range("H28").Value = 300
$ws.Range("K28").Value = 233.33

# row 29 - uy
$ws.Range("B29").Value = 325
$ws.Range("E29").Value = 325
$ws.Range("H29").Value = 200
$ws.Range("K29").Value = 325

# row 30 - sx
$ws.Range("B30").Value = 166.67
$ws.Range("E30").Value = 100
$ws.Range("H30").Value = 133.33
$ws.Range("K30").Value = 100

# row 31 - sy
$ws.Range("B31").Value = 75
$ws.Range("E31").Value = 50
$ws.Range("H31").Value = 100
$ws.Range("K31").Value = 50

# row 32 - theta (only the first column changes)
$ws.Range("B32").Value = 100

# row 33 - A
$ws.Range("B33").Value = 40
$ws.Range("E33").Value = 80
$ws.Range("H33").Value = 60
$ws.Range("K33").Value = 40

# --- Nudge the embedded surface chart a little (small drag in the Excel UI) ---
$co = $ws.ChartObjects().Item(1)
$co.Left = $co.Left - 1
$co.Top = $co.Top + 6

# --- Move the active selection to P26 ---
$ws.Range("P26").Select() | Out-Null
